{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Target edit (per the commit's diff):\n//   1. The paragraph \"Another git exercise\" is split into two runs \u2014\n//      \"Another git \" and \"exercise\" \u2014 with Word proofing-error markers\n//      (<w:proofErr w:type=\"gramStart\"/> ... <w:proofErr w:type=\"gramEnd\"/>)\n//      wrapped around the second run (the grammar checker flagging the\n//      sentence), and the first run gets xml:space=\"preserve\" because it\n//      now ends in a trailing space.\n//   2. A new, empty paragraph is added after it.\n//   3. A further new paragraph containing the text \"New line\" is added\n//      after that.\n//\n// `<w:proofErr>` markers aren't part of the Word JS object model (there's\n// no API that inserts them directly), so we build the exact OOXML for the\n// edited run content and splice it in with Range.insertOoxml \u2014 this is the\n// supported, documented escape hatch for OOXML the object model itself\n// can't author. insertOoxml requires the \"flat OPC\" <pkg:package> envelope.\n\nconst body = context.document.body;\n\n// --- Step 1: split \"Another git exercise\" into the two runs, with the\n// proofErr markers around \"exercise\", replacing just that text (not the\n// whole paragraph) so the paragraph's own identity/formatting is kept.\nconst splitOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">Another git </w:t></w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r><w:t>exercise</w:t></w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst targets = body.search(\"Another git exercise\", { matchCase: true });\nawait context.sync();\n\nif (targets.items.length === 0) {\n  throw new Error('Could not find \"Another git exercise\" in the document body.');\n}\n\ntargets.items[0].insertOoxml(splitOoxml, \"Replace\");\nawait context.sync();\n\n// --- Step 2: append a blank paragraph, then a \"New line\" paragraph, at\n// the very end of the body (after the paragraph we just edited).\nconst tailOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p/>\n          <w:p>\n            <w:r><w:t>New line</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst endRange = body.getRange(\"End\");\nendRange.insertOoxml(tailOoxml, \"After\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Target edit (per the commit's diff):\n#   1. The paragraph \"Another git exercise\" is split into two runs --\n#      \"Another git \" and \"exercise\" -- with Word's grammar-checker\n#      proofing-error markers (<w:proofErr w:type=\"gramStart\"/> ...\n#      <w:proofErr w:type=\"gramEnd\"/>) wrapped around the second run, and\n#      the first run gets xml:space=\"preserve\" because it now ends in a\n#      trailing space.\n#   2. A new, empty paragraph is added after it.\n#   3. A further new paragraph containing the text \"New line\" is added\n#      after that.\n#\n# <w:proofErr> markers have no dedicated property/method on the Word\n# object model (Word itself only ever inserts them as a side effect of\n# its grammar checker), so the supported way to author that exact markup\n# through COM is Range.InsertXML with a \"flat OPC\" <pkg:package> payload\n# -- the same mechanism backing Office.js's Range.insertOoxml.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: split \"Another git exercise\" into the two runs, with the\n# proofErr markers around \"exercise\". Find.Execute narrows $r1 down to\n# just that matched span, so InsertXML replaces only that text (not the\n# whole paragraph) and the paragraph itself keeps its own identity.\n$r1 = $d.Content\n$null = $r1.Find.Execute(\"Another git exercise\")\n\n$splitXml = @\"\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">Another git </w:t></w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r><w:t>exercise</w:t></w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$r1.InsertXML($splitXml)\n\n# --- Step 2: append a blank paragraph, then a \"New line\" paragraph, at\n# the very end of the document (after the paragraph we just edited).\n$endPos = $d.Content.End\n$r2 = $d.Range($endPos, $endPos)\n\n$tailXml = @\"\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p/>\n          <w:p>\n            <w:r><w:t>New line</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$r2.InsertXML($tailXml)\n"}
